$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("Bash", 0.72000000000000008, 0, 1, 0, 0.22499999999999912, 0.31499999999999911)
  ,@("C", 0.35999999999999915, 0, 1, 0, 0.45, 0.35999999999999915)
  ,@("C++", 0.40500000000000003, 0, 1, 0, 0.54, 0.31499999999999911)
  ,@("C#", 0.45, 0, 1, 0, 0.67500000000000004, 0.85499999999999998)
  ,@("Go", 0.63, 0, 0, 1, 0.45, 0.72000000000000008)
  ,@("HTML", 0.76500000000000001, 1, 0, 0, 0.76500000000000001, 0.54)
  ,@("Java", 0.54, 0, 1, 1, 0.81, 0.9)
  ,@("JavaScript", 0.85499999999999998, 1, 1, 0, 0.9, 0.63)
  ,@("JavaScript/React", 0.85499999999999998, 1, 0, 0, 0.9, 0.63)
  ,@("JavaScript/Node", 0.85499999999999998, 0, 1, 0, 0.9, 0.63)
  ,@("Kotlin", 0.58500000000000008, 0, 1, 0, 0.35999999999999915, 0.67500000000000004)
  ,@("MATLAB", 0.26999999999999907, 0, 0, 1, 0.22499999999999912, 0.26999999999999907)
  ,@("Perl", 0.17999999999999911, 0, 0, 1, 0.17999999999999911, 0.26999999999999907)
  ,@("PHP", 0.31499999999999911, 0, 1, 1, 0.63, 0.40500000000000003)
  ,@("Python/Flask", 0.9, 0, 0, 1, 0.7, 0.58500000000000008)
  ,@("Python/Django", 0.85, 0, 1, 1, 0.85, 0.58500000000000008)
  ,@("R", 0.22499999999999912, 0, 0, 1, 0.31499999999999911, 0.26999999999999907)
  ,@("Ruby", 0.81, 0, 1, 1, 0.49500000000000005, 0.45)
  ,@("Scala", 0.045000000000000005, 0, 1, 0, 0.26999999999999907, 0.49500000000000005)
  ,@("Scratch", 0.72000000000000008, 0, 0, 0, 0.13499999999999909, 0.26999999999999907)
  ,@("SQL", 0.13499999999999909, 0, 1, 1, 0.72000000000000008, 0.26999999999999907)
  ,@("Swift", 0.67500000000000004, 1, 1, 0, 0.40500000000000003, 0.76500000000000001)
  ,@("TypeScript", 0.49500000000000005, 1, 1, 1, 0.58500000000000008, 0.81)
  ,@("CSS", 0.77400000000000002, 1, 0, 0, 0.58500000000000008, 0.45)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $vals = $data[$i]
  $ws.Cells.Item($row, 1).Value = $vals[0]
  $ws.Cells.Item($row, 2).Value = $vals[1]
  $ws.Cells.Item($row, 3).Value = $vals[2]
  $ws.Cells.Item($row, 4).Value = $vals[3]
  $ws.Cells.Item($row, 5).Value = $vals[4]
  $ws.Cells.Item($row, 6).Value = $vals[5]
  $ws.Cells.Item($row, 7).Value = $vals[6]
}

$ws.Range("H14").Select()
